$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.930.97"
$ws.Range("E2").Value = "'  -3.25%  "

$ws.Range("D3").Value = "'2.911.59"
$ws.Range("E3").Value = "'  -3.67%  "

$ws.Range("E4").Value = "'  +0.00%  "

$ws.Range("D5").Value = "'591.04"
$ws.Range("E5").Value = "'  -0.85%  "

$ws.Range("D6").Value = "'145.05"
$ws.Range("E6").Value = "'  -5.21%  "

$ws.Range("E7").Value = "'  -0.06%  "

$ws.Range("E8").Value = "'  -1.57%  "

$ws.Range("D9").Value = "'2.909.08"
$ws.Range("E9").Value = "'  -3.56%  "

$ws.Range("E10").Value = "'  -3.95%  "

$ws.Range("E11").Value = "'  -3.63%  "

$ws.Range("D12").Value = "'0.444"
$ws.Range("E12").Value = "'  -4.14%  "

$ws.Range("E13").Value = "'  -2.55%  "

$ws.Range("D14").Value = "'33.60"
$ws.Range("E14").Value = "'  -5.90%  "

$ws.Range("D15").Value = "'0.127"
$ws.Range("E15").Value = "'  +0.29%  "

$ws.Range("D16").Value = "'3.396.16"
$ws.Range("E16").Value = "'  -3.62%  "

$ws.Range("D17").Value = "'60.881.67"
$ws.Range("E17").Value = "'  -3.27%  "

$ws.Range("D18").Value = "'6.74"
$ws.Range("E18").Value = "'  -4.75%  "

$ws.Range("D19").Value = "'2.905.70"
$ws.Range("E19").Value = "'  -3.86%  "

$ws.Range("D20").Value = "'430.57"
$ws.Range("E20").Value = "'  -4.17%  "

$ws.Range("D21").Value = "'13.55"
$ws.Range("E21").Value = "'  -4.75%  "

$ws.Range("E22").Value = "'  -1.90%  "

$ws.Range("E23").Value = "'  -5.73%  "

$ws.Range("D24").Value = "'81.35"
$ws.Range("E24").Value = "'  -1.84%  "

$ws.Range("D25").Value = "'10.83"
$ws.Range("E25").Value = "'  -4.34%  "

$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "'  -2.44%  "

$ws.Range("D27").Value = "'11.99"
$ws.Range("E27").Value = "'  -2.92%  "

$ws.Range("E28").Value = "'  +0.05%  "

$ws.Range("D29").Value = "'2.33"
$ws.Range("E29").Value = "'  +2.09%  "

$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "'  +0.10%  "

$ws.Range("D31").Value = "'2.63"
$ws.Range("E31").Value = "'  -2.59%  "

$ws.Range("D32").Value = "'7.08"
$ws.Range("E32").Value = "'  -5.65%  "

$ws.Range("D33").Value = "'26.67"
$ws.Range("E33").Value = "'  -3.46%  "

$ws.Range("E34").Value = "'  -2.97%  "

$ws.Range("D35").Value = "'0.0₃0856"
$ws.Range("E35").Value = "'  -1.86%  "

$ws.Range("E36").Value = "'  -2.96%  "

$ws.Range("E37").Value = "'  -4.40%  "

$ws.Range("D38").Value = "'3.02"
$ws.Range("E38").Value = "'  -3.47%  "

$ws.Range("D39").Value = "'49.63"
$ws.Range("E39").Value = "'  -1.65%  "

$ws.Range("E40").Value = "'  -3.49%  "

$ws.Range("D41").Value = "'2.01"
$ws.Range("E41").Value = "'  -4.70%  "

$ws.Range("E42").Value = "'  -4.42%  "

$ws.Range("D43").Value = "'0.293"
$ws.Range("E43").Value = "'  -3.68%  "

$ws.Range("D44").Value = "'40.30"
$ws.Range("E44").Value = "'  -9.97%  "

$ws.Range("B45").Value = "'Bittensor"
$ws.Range("C45").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'377.05"
$ws.Range("E45").Value = "'  -3.57%  "

$ws.Range("B46").Value = "'VeChain"
$ws.Range("C46").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0349"
$ws.Range("E46").Value = "'  -3.27%  "

$ws.Range("D47").Value = "'2.704.81"
$ws.Range("E47").Value = "'  +0.01%  "

$ws.Range("D48").Value = "'129.86"
$ws.Range("E48").Value = "'  -3.04%  "

$ws.Range("D50").Value = "'24.16"
$ws.Range("E50").Value = "'  -9.65%  "
